$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one row per calendar day (column A) with quantities in B:J.
# This edit rolls the date window forward by one day:
#   - drop the oldest day (2024-09-09, row 2) -> all rows shift up by one
#   - append the newest day (2024-10-08) as the new last row

$ws.Rows(2).Delete()

# Copy the formatting (font/border/alignment/style) of the last remaining
# date cell onto the new row so the appended row matches the existing
# column-A style (bold, bordered, centered) exactly.
$ws.Range("A29").Copy()
$ws.Range("A30").PasteSpecial(-4122)

# Write the new date as literal text (matching how every other date in
# column A is stored) rather than letting Excel auto-convert the
# date-like string into a serial date number: build it as a text formula
# result, then freeze it to a static value with Paste Values.
$ws.Range("A30").Formula = '="2024-10-08"'
$ws.Range("A30").Copy()
$ws.Range("A30").PasteSpecial(-4163)

# New day's quantities.
$ws.Range("B30").Value = 116.4121952
$ws.Range("C30").Value = 0.00170247
$ws.Range("D30").Value = 0.008850780000000001
$ws.Range("E30").Value = 0.06933635
$ws.Range("F30").Value = 12792.90181321
$ws.Range("G30").Value = 465.80531254
$ws.Range("H30").Value = 0.24
$ws.Range("I30").Value = 1.7904431
$ws.Range("J30").Value = 485.38834923

$excel.CutCopyMode = $false
